$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.882106333333333
$ws.Range("H2").Value = 11.646319
$ws.Range("I2").Value = 0.6257373677154582
$ws.Range("J2").Value = 0.6257373677154581
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.324764666666667
$ws.Range("N2").Value = 6.974294
$ws.Range("O2").Value = 0.04473923998638302
$ws.Range("P2").Value = 0.04473923998638301
$ws.Range("Q2").Value = 9.024983635976222
$ws.Range("R2").Value = 81.224852723786
$ws.Range("S2").Value = 0.02799501426266948
$ws.Range("T2").Value = 0.02799501426266947

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.882106333333333
$ws.Range("H3").Value = 11.646319
$ws.Range("I3").Value = 0.6257373677154582
$ws.Range("J3").Value = 0.6257373677154581
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 19.27491966666667
$ws.Range("N3").Value = 57.824759
$ws.Range("O3").Value = 0.3709387315842666
$ws.Range("P3").Value = 0.3709387315842665
$ws.Range("Q3").Value = 74.82728771245789
$ws.Range("R3").Value = 673.445589412121
$ws.Range("S3").Value = 0.2321102254852499
$ws.Range("T3").Value = 0.2321102254852498

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.882106333333333
$ws.Range("H4").Value = 11.646319
$ws.Range("I4").Value = 0.6257373677154582
$ws.Range("J4").Value = 0.6257373677154581
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 30.36285833333334
$ws.Range("N4").Value = 91.08857500000001
$ws.Range("O4").Value = 0.5843220284293504
$ws.Range("P4").Value = 0.5843220284293504
$ws.Range("Q4").Value = 117.8718446339361
$ws.Range("R4").Value = 1060.846601705425
$ws.Range("S4").Value = 0.3656321279675389
$ws.Range("T4").Value = 0.3656321279675388

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.522503666666667
$ws.Range("H5").Value = 4.567511000000001
$ws.Range("I5").Value = 0.2454047764062963
$ws.Range("J5").Value = 0.2454047764062963
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.324764666666667
$ws.Range("N5").Value = 6.974294
$ws.Range("O5").Value = 0.04473923998638302
$ws.Range("P5").Value = 0.04473923998638301
$ws.Range("Q5").Value = 3.539462729137112
$ws.Range("R5").Value = 31.85516456223401
$ws.Range("S5").Value = 0.01097922318544596
$ws.Range("T5").Value = 0.01097922318544595

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.522503666666667
$ws.Range("H6").Value = 4.567511000000001
$ws.Range("I6").Value = 0.2454047764062963
$ws.Range("J6").Value = 0.2454047764062963
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 19.27491966666667
$ws.Range("N6").Value = 57.824759
$ws.Range("O6").Value = 0.3709387315842666
$ws.Range("P6").Value = 0.3709387315842665
$ws.Range("Q6").Value = 29.34613586720545
$ws.Range("R6").Value = 264.1152228048491
$ws.Range("S6").Value = 0.09103013648487211
$ws.Range("T6").Value = 0.09103013648487208

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.522503666666667
$ws.Range("H7").Value = 4.567511000000001
$ws.Range("I7").Value = 0.2454047764062963
$ws.Range("J7").Value = 0.2454047764062963
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 30.36285833333334
$ws.Range("N7").Value = 91.08857500000001
$ws.Range("O7").Value = 0.5843220284293504
$ws.Range("P7").Value = 0.5843220284293504
$ws.Range("Q7").Value = 46.22756314298056
$ws.Range("R7").Value = 416.0480682868251
$ws.Range("S7").Value = 0.1433954167359783
$ws.Range("T7").Value = 0.1433954167359782

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 0.7994406666666668
$ws.Range("H8").Value = 2.398322
$ws.Range("I8").Value = 0.1288578558782456
$ws.Range("J8").Value = 0.1288578558782456
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.324764666666667
$ws.Range("N8").Value = 6.974294
$ws.Range("O8").Value = 0.04473923998638302
$ws.Range("P8").Value = 0.04473923998638301
$ws.Range("Q8").Value = 1.858511414963111
$ws.Range("R8").Value = 16.726602734668
$ws.Range("S8").Value = 0.005765002538267586
$ws.Range("T8").Value = 0.005765002538267584

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 0.7994406666666668
$ws.Range("H9").Value = 2.398322
$ws.Range("I9").Value = 0.1288578558782456
$ws.Range("J9").Value = 0.1288578558782456
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 19.27491966666667
$ws.Range("N9").Value = 57.824759
$ws.Range("O9").Value = 0.3709387315842666
$ws.Range("P9").Value = 0.3709387315842665
$ws.Range("Q9").Value = 15.40915462826645
$ws.Range("R9").Value = 138.682391654398
$ws.Range("S9").Value = 0.04779836961414466
$ws.Range("T9").Value = 0.04779836961414464

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.7994406666666668
$ws.Range("H10").Value = 2.398322
$ws.Range("I10").Value = 0.1288578558782456
$ws.Range("J10").Value = 0.1288578558782456
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 30.36285833333334
$ws.Range("N10").Value = 91.08857500000001
$ws.Range("O10").Value = 0.5843220284293504
$ws.Range("P10").Value = 0.5843220284293504
$ws.Range("Q10").Value = 24.27330370790556
$ws.Range("R10").Value = 218.45973337115
$ws.Range("S10").Value = 0.07529448372583337
$ws.Range("T10").Value = 0.07529448372583336
